# Fruta / hortaliza, semanal
# Inserts one new weekly price record as row 697 in the "Uva" sheet,
# pushing the existing rows 697-775 down to 698-776.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 697 (shifts 697..775 -> 698..776)
$ws.Rows.Item(697).Insert()

# Populate the new row 697 with the new weekly record
$ws.Range("A697").Value = 5
$ws.Range("B697").Value = "Macroferia Regional de Talca"
$ws.Range("C697").Value = "Maule"
$ws.Range("D697").Value = 45124
$ws.Range("E697").Value = 7
$ws.Range("F697").Value = "Fruta"
$ws.Range("G697").Value = 100109
$ws.Range("H697").Value = "Uva"
$ws.Range("I697").Value = 100109001
$ws.Range("J697").Value = "Uva"
$ws.Range("K697").Value = "Crimpson Seedless"
$ws.Range("L697").Value = "Primera"
$ws.Range("M697").Value = 100
$ws.Range("N697").Value = 12000
$ws.Range("O697").Value = 12000
$ws.Range("P697").Value = 12000
$ws.Range("Q697").Value = "`$/bandeja 8 kilos"
$ws.Range("R697").Value = "Región de O'Higgins"
$ws.Range("S697").Value = 1500
$ws.Range("T697").Value = 8
